$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-06-10T12:41:12+00:00"

# --- Concepts sheet: fix capitalization of several Display values ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C5").Value = "Muscle Biopsy"
$concepts.Range("C6").Value = "Metabolic Work-Up"
$concepts.Range("C7").Value = "Serum Creatine Kinase"
$concepts.Range("C8").Value = "Plasma Amino Acid Chromatography"
$concepts.Range("C21").Value = "GCN Repeat Testing (Oculopharyngeal Muscular Dystrophy)"
$concepts.Range("C22").Value = "Deletions and Duplications Testing (Duchenne and Becker Dystrophies)"

$wb.Save()
